# Weekly refresh of Fruta / hortaliza data.
# Updates Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) for rows 2-9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 45092; J = 210; K = 10000; L = 11000; M = 10714; P = 595 }
    3 = @{ D = 45245; J = 100; K = 9000;  L = 10000; M = 9500;  P = 528 }
    4 = @{ D = 44792; J = 160; K = 9000;  L = 10000; M = 9500;  P = 528 }
    5 = @{ D = 45205; J = 200; K = 11000; L = 12000; M = 11500; P = 639 }
    6 = @{ D = 44804; J = 50;  K = 9500;  L = 10000; M = 9750;  P = 542 }
    7 = @{ D = 45215; J = 200; K = 11000; L = 12000; M = 11500; P = 639 }
    8 = @{ D = 45175; J = 250; K = 11000; L = 12000; M = 11500; P = 639 }
    9 = @{ D = 44714; J = 80;  K = 9000;  L = 10000; M = 9500;  P = 528 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
